# Bill_Statement.xlsx: populate item rows 2-22 on the active sheet,
# filling in item name, quantity, price and tax for each menu item.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Write-Host "Starting Bill_Statement data load..."

$items = @(
    @("Vanilla Cream", 50, 8750, 700),
    @("Filter Coffee", 30, 2700, 216),
    @("Cold Coffee", 12, 1800, 144),
    @("Black Coffee", 10, 880, 70.40000000000001),
    @("Normal Tea", 10, 250, 20),
    @("Lemon Tea", 10, 350, 28),
    @("Ginger Tea", 10, 400, 32),
    @("Honey Tea", 10, 405, 32.4),
    @("Bhel Puri", 10, 2200, 264),
    @("Sev Puri", 10, 1700, 204),
    @("Peanut Masala", 10, 2990, 358.8),
    @("Chilli Cheese", 10, 3450, 414),
    @("Chicken Tikki", 10, 2980, 536.4),
    @("Kheema chat", 10, 3100, 558),
    @("Shambi Kabab", 10, 3750, 675),
    @("Chicken Tandoor", 10, 3400, 612),
    @("Veggie Delite", 10, 2500, 300),
    @("Paneer Tikka", 10, 2750, 330),
    @("Corn n Peas", 10, 1750, 210),
    @("Peri Peri Chicken", 10, 3350, 603),
    @("Ham Chicken", 10, 3250, 585),
)

$row = 2
foreach ($item in $items) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    Write-Host "Wrote row $row : $($item[0])"
    $row++
}

Write-Host "Bill_Statement data load complete: wrote $($items.Count) rows."
